$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying "2022-Q3" (same
#    fund-holding layout/formatting), inserting it right after "总计"
#    and before "2022-Q3", then updating the four changed data points.
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q3")
$src.Copy($wb.Worksheets.Item(2))
$newQ4 = $wb.Worksheets.Item(2)
$newQ4.Name = "2022-Q4"

$newQ4.Range("D2").Value = "'1.45"
$newQ4.Range("D2").Style = "Normal"
$newQ4.Range("E2").Value = "'78.91"
$newQ4.Range("E2").Style = "Normal"
$newQ4.Range("F2").Value = "'5.75"
$newQ4.Range("F2").Style = "Normal"
$newQ4.Range("G2").Value = "'0.0834"
$newQ4.Range("G2").Style = "Normal"

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row right below the
#    header for 2022-Q4 and push the existing quarters down one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows(2).Insert()

# Give the new A2 the same bold/bordered look as the rest of column A
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").Style = "Normal"

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.08

$summary.Range("A3").Value = 1
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.09

$summary.Range("A4").Value = 2
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.13

$summary.Range("A5").Value = 3
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 0.25

$summary.Range("A6").Value = 4
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 0.08

$summary.Range("A7").Value = 5
$summary.Range("C7").Value = 12
$summary.Range("D7").Value = 2.27

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q2"
$summary.Range("C8").Value = 6
$summary.Range("D8").Value = 0.84

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2021-Q1"
$summary.Range("C9").Value = 3
$summary.Range("D9").Value = 0.34
